$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 40
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 20
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 30

$ws.Range("B7").Select()
